# Add a new "14-jun" column (J) to the sheet, mirroring the existing
# date columns (C..I) both in header text and in per-row values/format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell, same style as the other date headers (e.g. C1/D1).
$ws.Range("J1").Value = "14-jun"
$ws.Range("J1").NumberFormat = "@"

# New per-row numeric values, copying the number format/style used by
# the existing "13-jun" column (I).
$values = @(
    0,
    13.132209895812792,
    17.626869978956577,
    19.421960383193653,
    0,
    7.1148367049590693,
    5.4360915226953104,
    10.979997668130585,
    14.115669720462879,
    15.279211001908271,
    0,
    14.594242540544988,
    0,
    0,
    12.510503873694049,
    0,
    0
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 10)
    $cell.Value = $values[$i]
    $cell.NumberFormat = "0"
}

# Mirror the selection change captured in the edit (J2:J18 selected,
# anchored at J2).
$ws.Range("J2:J18").Select()
